# Add "R" to the "Experienced" skills bullet in the Languages & Technologies
# section, e.g. "...(several paradigms), Python " -> "...Python, R ".
# Word also re-homes its internal "_GoBack" bookmark (last-edit marker) from
# wherever it previously sat (after "Familiar" in the next bullet) to right
# after this newly typed text, which is what a live edit in Word itself
# would do.

$d = $word.ActiveDocument

# Locate the "Experienced: ..." bullet paragraph (the skills line that ends
# in "... Python ") rather than relying on a hard-coded paragraph index.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "^Experienced:.*Python") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    # Find "Python" inside that paragraph, collapse to just after it, then
    # grow the collapsed point by one character so it spans the single
    # trailing space run that follows "Python". Rewriting that run's text
    # (instead of blindly inserting at the boundary) keeps the preceding
    # ", " run untouched and simply extends "Python" into "Python, R ",
    # followed by a fresh trailing-space run - matching how Word itself
    # lays the runs out after a real edit in this spot.
    $f = $target.Range.Duplicate
    $f.Find.Execute("Python", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $f.Collapse(0)
    $f.MoveEnd(1, 1) | Out-Null
    $f.Text = ", R "

    # Move the "_GoBack" bookmark (Word's marker for the last edited spot)
    # to sit right after the text we just typed.
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }

    $bmRange = $target.Range.Duplicate
    $bmRange.Find.Execute(", R", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $bmRange.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
